$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 13: First Name / Mid Name / Last Name / Concatinate ---
# A13 already holds "Name" (header style s=1). Match the bold/centered
# header style for B13:D13, and the centered (non-bold) style for E13.
$ws.Range("B13").Font.Bold = $true
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("B13").VerticalAlignment = -4108
$ws.Range("B13").Value = "First Name"

$ws.Range("C13").Font.Bold = $true
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("C13").VerticalAlignment = -4108
$ws.Range("C13").Value = "Mid Name"

$ws.Range("D13").Font.Bold = $true
$ws.Range("D13").HorizontalAlignment = -4108
$ws.Range("D13").VerticalAlignment = -4108
$ws.Range("D13").Value = "Last Name"

$ws.Range("E13").HorizontalAlignment = -4108
$ws.Range("E13").VerticalAlignment = -4108
$ws.Range("E13").Value = "Concatinate"

# --- Row 14: Shivendra Singh Mira ---
$ws.Range("B14").Formula = '=LEFT(A14,FIND(" ",A14) - 1)'
$ws.Range("C14").Formula = '=MID(A15,FIND(" ",A15)+1,FIND(" ",A15,FIND(" ",A15) +1-FIND(" ",A15)))'
$ws.Range("D14").Formula = '=RIGHT(A14, LEN(A14) - FIND(" ", A14, FIND(" ",A14) +1))'
$ws.Range("E14").Formula = '=TRIM(CONCATENATE(B14," ", C14," ",D14))'

# --- Row 15 (and the filled-down row 16): First-name / Concatinate columns are
# entered as one relative fill across B15:B16 and E15:E16 respectively, which
# is how they end up sharing a single formula definition (same as the source
# workbook: si="11" for LEFT(...), si="12" for TRIM(CONCATENATE(...))).
$ws.Range("B15:B16").Formula = '=LEFT(A15,FIND(" ",A15) - 1)'
$ws.Range("C15").Formula = '=SUBSTITUTE(A15,B14,"")'
$ws.Range("D15").Formula = '=IFERROR(RIGHT(A15, LEN(A15) - FIND(" ", A15, FIND(" ",A15) +1)),"")'
$ws.Range("E15:E16").Formula = '=TRIM(CONCATENATE(B15," ", C15," ",D15))'

# --- Row 16: Noopur Sidhartha Vishwakarma (B16/E16 already filled above) ---
$ws.Range("D16").Formula = '=RIGHT(A16, LEN(A16) - FIND(" ", A16, FIND(" ",A16) +1))'

# --- Window / selection bookkeeping ---
$ws.Range("C15").Select() | Out-Null
